$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.525288526816022
$ws.Range("F2").Value = 0.945350734094617
$ws.Range("J2").Value = 0.0072463768115942
$ws.Range("K2").Value = 0.193981481481481
$ws.Range("L2").Value = 0.963532356095817
$ws.Range("M2").Value = 0.845771144278607
$ws.Range("P2").Value = 38.5842981036431
$ws.Range("Q2").Value = 84.2261089565554
$ws.Range("R2").Value = 36.4842671200365
$ws.Range("S2").Value = 36.444452276668
$ws.Range("T2").Value = 94.68480666949
$ws.Range("U2").Value = 4.70990948609626
$ws.Range("V2").Value = 11.6424298078935
$ws.Range("W2").Value = 4.23926165985295
$ws.Range("X2").Value = 4.06213030060682
$ws.Range("Y2").Value = 11.1603131756513
$ws.Range("Z2").Value = 21183.1629166667
$ws.Range("AB2").Value = 387.33756130081
$ws.Range("AC2").Value = 13496.1731944444
$ws.Range("AD2").Value = 7299.65216092141
$ws.Range("AE2").Value = 2300.34783907859
$ws.Range("F3").Value = 0.210032626427406
$ws.Range("G3").Value = 0.507142857142857
$ws.Range("H3").Value = 0.999680715197957
$ws.Range("L3").Value = 0.289473684210526
$ws.Range("M3").Value = 0.739287558079504
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 34.4541272545617
$ws.Range("Q3").Value = 73.577675587252
$ws.Range("R3").Value = 33.0308755311829
$ws.Range("S3").Value = 31.8073515075635
$ws.Range("T3").Value = 84.3299492697222
$ws.Range("U3").Value = 3.41241679533663
$ws.Range("V3").Value = 9.60524408162853
$ws.Range("W3").Value = 3.57101220662391
$ws.Range("X3").Value = 3.45252777291113
$ws.Range("Y3").Value = 8.98723844336624
$ws.Range("Z3").Value = 21183.1629166667
$ws.Range("AB3").Value = 405.40939274852
$ws.Range("AC3").Value = 7009.50152777778
$ws.Range("AD3").Value = 13768.2519961404
$ws.Range("AE3").Value = 5431.74800385963
$ws.Range("E4").Value = 0
$ws.Range("H4").Value = 0.752554278416347
$ws.Range("L4").Value = 0.0551595383570944
$ws.Range("M4").Value = 0.181694546443628
$ws.Range("O4").Value = 0.292903225806452
$ws.Range("P4").Value = 33.2648156459609
$ws.Range("Q4").Value = 70.8879611879411
$ws.Range("R4").Value = 32.3006923372185
$ws.Range("S4").Value = 30.887821394275
$ws.Range("T4").Value = 82.2171261917537
$ws.Range("U4").Value = 3.19692489968859
$ws.Range("V4").Value = 9.26414868294762
$ws.Range("W4").Value = 3.48973551488594
$ws.Range("X4").Value = 3.36294078389247
$ws.Range("Y4").Value = 8.78865054319276
$ws.Range("Z4").Value = 21183.1629166667
$ws.Range("AB4").Value = 314.697247910845
$ws.Range("AC4").Value = 3119.25388888889
$ws.Range("AD4").Value = 17749.2117798669
$ws.Range("AE4").Value = 11050.7882201331
$ws.Range("F5").Value = 0.000135943447525829
$ws.Range("H5").Value = 0.127394636015326
$ws.Range("L5").Value = 0.00305498981670061
$ws.Range("M5").Value = 0.0375254928619986
$ws.Range("O5").Value = 0.770581778265642
$ws.Range("P5").Value = 33.0792760551785
$ws.Range("Q5").Value = 70.6285069576724
$ws.Range("R5").Value = 32.2992252226806
$ws.Range("S5").Value = 30.8346113660416
$ws.Range("T5").Value = 81.9768821478371
$ws.Range("U5").Value = 3.17759358972222
$ws.Range("V5").Value = 9.25016903961268
$ws.Range("W5").Value = 3.47923476930583
$ws.Range("X5").Value = 3.34190152189813
$ws.Range("Y5").Value = 8.79643039198368
$ws.Range("Z5").Value = 21183.1629166667
$ws.Range("AB5").Value = 298.736110537288
$ws.Range("AC5").Value = 527.342638888889
$ws.Range("AD5").Value = 20357.0841672405
$ws.Range("AE5").Value = 18042.9158327595
$ws.Range("H6").Value = 0
$ws.Range("M6").Value = 0.0067990209409845
$ws.Range("O6").Value = 0.278416347381865
$ws.Range("P6").Value = 32.6264185247354
$ws.Range("Q6").Value = 69.5646243247223
$ws.Range("R6").Value = 31.6679917674393
$ws.Range("S6").Value = 30.2792756865068
$ws.Range("T6").Value = 80.5455115666195
$ws.Range("U6").Value = 3.09129012886463
$ws.Range("V6").Value = 9.07331281530734
$ws.Range("W6").Value = 3.43168568904769
$ws.Range("X6").Value = 3.25957677391981
$ws.Range("Y6").Value = 8.66223381950367
$ws.Range("Z6").Value = 21183.1629166667
$ws.Range("AB6").Value = 295.356514468967
$ws.Range("AC6").Value = 2.60375
$ws.Range("AD6").Value = 20885.2026521977
$ws.Range("AE6").Value = 27114.7973478023
$ws.Range("H7").Value = 0
$ws.Range("M7").Value = 0.00013598041881969
$ws.Range("O7").Value = 0.055874840357599
$ws.Range("P7").Value = 32.4246142298623
$ws.Range("Q7").Value = 69.0646560599087
$ws.Range("R7").Value = 31.3984739617419
$ws.Range("S7").Value = 29.8694418645876
$ws.Range("T7").Value = 79.9697483038763
$ws.Range("U7").Value = 3.0603735754767
$ws.Range("V7").Value = 8.95035178117271
$ws.Range("W7").Value = 3.40712230876119
$ws.Range("X7").Value = 3.24637844786093
$ws.Range("Y7").Value = 8.59499915469603
$ws.Range("Z7").Value = 21183.1629166667
$ws.Range("AB7").Value = 294.226518834546
$ws.Range("AC7").Value = 2.60375
$ws.Range("AD7").Value = 20886.3326478321
$ws.Range("AE7").Value = 36713.6673521679
$ws.Range("F8").Value = 0.000271886895051658
$ws.Range("P8").Value = 32.36706652199
$ws.Range("Q8").Value = 68.9532748682699
$ws.Range("R8").Value = 31.368157773898
$ws.Range("S8").Value = 29.783714812611
$ws.Range("T8").Value = 79.7346847524877
$ws.Range("U8").Value = 3.06009579769892
$ws.Range("V8").Value = 8.9424451681594
$ws.Range("W8").Value = 3.39880273118644
$ws.Range("X8").Value = 3.24637844786093
$ws.Range("Y8").Value = 8.53779072713542
$ws.Range("Z8").Value = 21183.1629166667
$ws.Range("AB8").Value = 293.952537168908
$ws.Range("AC8").Value = 2.60375
$ws.Range("AD8").Value = 20886.6066294978
$ws.Range("AE8").Value = 46313.3933705022
$ws.Range("F9").Value = 0.000271886895051658
$ws.Range("P9").Value = 32.360993478445
$ws.Range("Q9").Value = 68.9051064625343
$ws.Range("R9").Value = 31.3396421410115
$ws.Range("S9").Value = 29.7670336700955
$ws.Range("T9").Value = 79.705429937693
$ws.Range("U9").Value = 3.05105333030138
$ws.Range("V9").Value = 8.93282940238792
$ws.Range("W9").Value = 3.39880273118644
$ws.Range("X9").Value = 3.24637844786093
$ws.Range("Y9").Value = 8.53182369470614
$ws.Range("Z9").Value = 21183.1629166667
$ws.Range("AB9").Value = 293.824481613353
$ws.Range("AC9").Value = 2.60375
$ws.Range("AD9").Value = 20886.7346850533
$ws.Range("AE9").Value = 55913.2653149467
